# Auto-generated edit script: updates crypto price/volume table to match latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.674.15"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.689.26"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'315.41"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.3941"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'1.487"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").Value = "'1.002"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'0.08832"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'7.232"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'23.50"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'8.031"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").Value = "'0.00001315"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "1.694.56"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'99.55"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "'0.07016"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'6.984"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'14.30"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "24.654.25"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'3.305"
$ws.Range("E25").Value = "  +9.89%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "'162.41"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'135.25"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "'5.173"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").Value = "'7.630"
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("D32").Value = "1.878.65"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "'1.058"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "'0.08539"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "'7.093"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").Value = "'11.24"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "'0.2731"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'1.888"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").Value = "'0.09190"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "'0.02717"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "'0.7616"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "'16.00"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").Value = "'2.585"
$ws.Range("E45").Value = "  +4.67%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'4.214"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "'139.71"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "'0.07969"
$ws.Range("E51").Value = "  -0.57%  "
